# Update statistics for 10 DEZ -- add new review period column K ("19 NOV - 09 DEZ")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

# --- New column K: period header / date ---
$ws.Range("K1").Value = 41982
$ws.Range("K1").NumberFormat = "dd/mm/yyyy"

$ws.Range("K2").Value = "19 NOV - 09 DEZ"
$ws.Range("K2").NumberFormat = "0.00"

# --- bugs / features / infra block (rows 4-5) ---
$ws.Range("K4").Value = 15
$ws.Range("K4").NumberFormat = "0.00"

$ws.Range("K5").Value = 48
$ws.Range("K5").NumberFormat = "0.00"

# --- number of tickets (rows 9-12) ---
$ws.Range("K9").Value = 8
$ws.Range("K9").NumberFormat = "0.00"

$ws.Range("K10").Value = 17
$ws.Range("K10").NumberFormat = "0.00"

$ws.Range("K11").Value = 7
$ws.Range("K11").NumberFormat = "0.00"

$ws.Range("K12").Formula = "=SUM(K9:K11)"
$ws.Range("K12").NumberFormat = "0.00"

# --- tracked days on tickets (rows 14-17) ---
$ws.Range("K14").Value = 8
$ws.Range("K14").NumberFormat = "0.00"

$ws.Range("K15").Value = 32
$ws.Range("K15").NumberFormat = "0.00"

$ws.Range("K16").Value = 2
$ws.Range("K16").NumberFormat = "0.00"

$ws.Range("K17").Formula = "=SUM(K14:K16)"
$ws.Range("K17").NumberFormat = "0.00"

# --- tickets > 4 / > 1 tracked days (rows 18-19) ---
$ws.Range("K18").Value = 12
$ws.Range("K18").NumberFormat = "0.00"

$ws.Range("K19").Value = 3
$ws.Range("K19").NumberFormat = "0.00"

# --- row 20 stays blank (only spans metadata changes, handled by engine) ---

# --- miscalculated estimations block (rows 21, 25) ---
$ws.Range("K21").Value = 3
$ws.Range("K21").NumberFormat = "0.00"

$ws.Range("K25").Value = 1.1000000000000001
$ws.Range("K25").NumberFormat = "0.00"
$ws.Range("C25").Formula = "=AVERAGE(D25:K25)"

# --- bug issues open/closed (rows 27-29) ---
$ws.Range("K27").Value = 40
$ws.Range("K27").NumberFormat = "0.00"

$ws.Range("K28").Value = 298
$ws.Range("K28").NumberFormat = "0.00"

$ws.Range("K29").Formula = "=SUM(K27:K28)"
$ws.Range("K29").NumberFormat = "0.00"

# --- EVALUATION block (rows 31-34) ---
$ws.Range("K31").Formula = "=(K18/K12)"
$ws.Range("K31").NumberFormat = "0.00"
$ws.Range("C31").Formula = "=AVERAGE(D31:K31)"

$ws.Range("K32").Formula = "=(K19/K12)"
$ws.Range("K32").NumberFormat = "0.00"

$ws.Range("K33").Formula = "=K17/K5"
$ws.Range("K33").NumberFormat = "0.00"
$ws.Range("C33").Formula = "=AVERAGE(D33:K33)"

$ws.Range("K34").Formula = "=K12/K5"
$ws.Range("K34").NumberFormat = "0.00"
$ws.Range("C34").Formula = "=AVERAGE(D34:K34)"

# --- MITTELWERTE percent block (rows 36-38) ---
$ws.Range("K36").Formula = "=(K9/K12*100)"
$ws.Range("K36").NumberFormat = "0.00"

$ws.Range("K37").Formula = "=(K10/K12*100)"
$ws.Range("K37").NumberFormat = "0.00"

$ws.Range("K38").Formula = "=(K11/K12*100)"
$ws.Range("K38").NumberFormat = "0.00"

Write-Host "done"
